# Weekly refresh: push a new week's worth of data (2 rows) onto the top of
# the historical block (rows 232:345) by inserting 2 rows at row 232 -
# shifting the existing history down to 234:347 - and then filling the
# freshly inserted rows 232:233 with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 232; everything that was on rows
# 232:345 moves down to 234:347 (formats/styles travel with it).
$ws.Rows("232:233").Insert()

# Row 232 - "Primera" quality entry for the new date.
$ws.Cells.Item(232, 1).Value = 8
$ws.Cells.Item(232, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(232, 3).Value = "Coquimbo"
$ws.Cells.Item(232, 4).Value = 44845
$ws.Cells.Item(232, 5).Value = 4
$ws.Cells.Item(232, 6).Value = 100114014
$ws.Cells.Item(232, 7).Value = "Betarraga"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 2400
$ws.Cells.Item(232, 11).Value = 550
$ws.Cells.Item(232, 12).Value = 600
$ws.Cells.Item(232, 13).Value = 575
$ws.Cells.Item(232, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(232, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(232, 16).Value = 192
$ws.Cells.Item(232, 17).Value = 3
$ws.Cells.Item(232, 18).Value = "Hortaliza"

# Row 233 - "Segunda" quality entry for the new date.
$ws.Cells.Item(233, 1).Value = 8
$ws.Cells.Item(233, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(233, 3).Value = "Coquimbo"
$ws.Cells.Item(233, 4).Value = 44845
$ws.Cells.Item(233, 5).Value = 4
$ws.Cells.Item(233, 6).Value = 100114014
$ws.Cells.Item(233, 7).Value = "Betarraga"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Segunda"
$ws.Cells.Item(233, 10).Value = 1560
$ws.Cells.Item(233, 11).Value = 450
$ws.Cells.Item(233, 12).Value = 500
$ws.Cells.Item(233, 13).Value = 475
$ws.Cells.Item(233, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(233, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(233, 16).Value = 158
$ws.Cells.Item(233, 17).Value = 3
$ws.Cells.Item(233, 18).Value = "Hortaliza"
